$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the at-risk cells as Text so Excel keeps the literal
# string instead of silently converting it to a floating point number
# (which would drop trailing zeros / alter precision).
$textCells = @("D5", "D6", "D8", "D12", "D15", "D20", "D21", "D22", "D23", "D24", "D28", "D30", "D33", "D34", "D35", "D37", "D39", "D40", "D43", "D44", "D46", "D47", "D48", "D49", "D50")
foreach ($addr in $textCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "64.808.63"
$ws.Range("E2").Value = "  -1.87%  "
$ws.Range("D3").Value = "3.125.83"
$ws.Range("E3").Value = "  -7.84%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "568.51"
$ws.Range("D6").Value = "168.49"
$ws.Range("E6").Value = "  -6.05%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "0.601"
$ws.Range("E8").Value = "  -3.56%  "
$ws.Range("D9").Value = "3.124.38"
$ws.Range("E9").Value = "  -7.94%  "
$ws.Range("E10").Value = "  -5.91%  "
$ws.Range("E11").Value = "  -5.94%  "
$ws.Range("D12").Value = "0.388"
$ws.Range("E12").Value = "  -5.83%  "
$ws.Range("D13").Value = "3.663.73"
$ws.Range("E13").Value = "  -8.07%  "
$ws.Range("E14").Value = "  +1.04%  "
$ws.Range("D15").Value = "26.54"
$ws.Range("E15").Value = "  -8.44%  "
$ws.Range("D16").Value = "64.701.28"
$ws.Range("E16").Value = "  -2.18%  "
$ws.Range("E17").Value = "  -6.32%  "
$ws.Range("D18").Value = "3.128.44"
$ws.Range("E18").Value = "  -8.58%  "
$ws.Range("E19").Value = "  -3.17%  "
$ws.Range("D20").Value = "12.66"
$ws.Range("E20").Value = "  -7.25%  "
$ws.Range("D21").Value = "353.66"
$ws.Range("E21").Value = "  -3.37%  "
$ws.Range("D22").Value = "7.19"
$ws.Range("E22").Value = "  -4.36%  "
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.24%  "
$ws.Range("D24").Value = "68.73"
$ws.Range("E24").Value = "  -5.63%  "
$ws.Range("D26").Value = "3.268.66"
$ws.Range("E26").Value = "  -7.82%  "
$ws.Range("E27").Value = "  -8.13%  "
$ws.Range("D28").Value = "9.56"
$ws.Range("E28").Value = "  -1.73%  "
$ws.Range("E29").Value = "  -2.48%  "
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("E32").Value = "  -4.15%  "
$ws.Range("D33").Value = "21.69"
$ws.Range("E33").Value = "  -6.29%  "
$ws.Range("D34").Value = "5.21"
$ws.Range("E34").Value = "  -9.05%  "
$ws.Range("D35").Value = "6.54"
$ws.Range("E35").Value = "  -6.28%  "
$ws.Range("E36").Value = "  -5.46%  "
$ws.Range("D37").Value = "158.05"
$ws.Range("E37").Value = "  -2.20%  "
$ws.Range("E38").Value = "  -6.64%  "
$ws.Range("D39").Value = "0.824"
$ws.Range("E39").Value = "  -3.78%  "
$ws.Range("D40").Value = "26.15"
$ws.Range("E40").Value = "  -3.43%  "
$ws.Range("E41").Value = "  -1.67%  "
$ws.Range("D42").Value = "2.634.93"
$ws.Range("E42").Value = "  -1.68%  "
$ws.Range("D43").Value = "6.06"
$ws.Range("E43").Value = "  -3.10%  "
$ws.Range("D44").Value = "2.38"
$ws.Range("E44").Value = "  -8.17%  "
$ws.Range("E45").Value = "  -4.50%  "
$ws.Range("D46").Value = "39.31"
$ws.Range("E46").Value = "  -0.91%  "
$ws.Range("D47").Value = "0.0648"
$ws.Range("E47").Value = "  -4.10%  "
$ws.Range("D48").Value = "23.70"
$ws.Range("E48").Value = "  -3.05%  "
$ws.Range("D49").Value = "316.43"
$ws.Range("E49").Value = "  -5.58%  "
$ws.Range("D50").Value = "0.0269"
$ws.Range("E50").Value = "  -5.22%  "
$ws.Range("E51").Value = "  -2.00%  "
